$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.706.60'
$ws.Range('E2').Value = '  +0.43%  '
$ws.Range('D3').Value = '3.321.73'
$ws.Range('E3').Value = '  +1.48%  '
$ws.Range('E4').Value = '  +0.75%  '
$ws.Range('D5').Value = '''519.08'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -0.10%  '
$ws.Range('D6').Value = '''171.32'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -5.09%  '
$ws.Range('D7').Value = '''0.588'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  -2.61%  '
$ws.Range('D8').Value = '3.325.89'
$ws.Range('E8').Value = '  +2.03%  '
$ws.Range('E9').Value = '  +0.16%  '
$ws.Range('D10').Value = '''0.601'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -2.94%  '
$ws.Range('D11').Value = '''52.60'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -10.86%  '
$ws.Range('D12').Value = '''0.132'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  -0.78%  '
$ws.Range('D13').Value = '''0.0000254'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  -1.51%  '
$ws.Range('D14').Value = '''8.92'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -2.15%  '
$ws.Range('D15').Value = '3.892.03'
$ws.Range('E15').Value = '  +3.17%  '
$ws.Range('D16').Value = '3.348.41'
$ws.Range('E16').Value = '  +2.94%  '
$ws.Range('E17').Value = '  -0.90%  '
$ws.Range('B18').Value = 'WrappedBTC'
$ws.Range('C18').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D18').Value = '63.647.72'
$ws.Range('E18').Value = '  +1.00%  '
$ws.Range('B19').Value = 'Chainlink'
$ws.Range('C19').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D19').Value = '''17.37'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -1.23%  '
$ws.Range('D20').Value = '''11.11'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +0.57%  '
$ws.Range('D21').Value = '''0.950'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -0.70%  '
$ws.Range('D22').Value = '''371.23'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -0.19%  '
$ws.Range('D23').Value = '''4.23'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +8.54%  '
$ws.Range('D24').Value = '''11.29'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +0.78%  '
$ws.Range('D25').Value = '''81.03'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +0.77%  '
$ws.Range('D26').Value = '''3.64'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -1.92%  '
$ws.Range('D27').Value = '''6.18'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +1.93%  '
$ws.Range('D28').Value = '''2.67'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +0.57%  '
$ws.Range('D29').Value = '''11.14'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -3.31%  '
$ws.Range('D30').Value = '''8.10'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -3.42%  '
$ws.Range('D31').Value = '''28.54'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -0.58%  '
$ws.Range('D32').Value = '''624.93'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -1.62%  '
$ws.Range('D33').Value = '''6.35'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -7.71%  '
$ws.Range('D34').Value = '''11.09'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -2.19%  '
$ws.Range('D35').Value = '''0.104'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -1.25%  '
$ws.Range('D36').Value = '''57.60'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -1.66%  '
$ws.Range('D37').Value = '''0.999'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -0.01%  '
$ws.Range('D38').Value = '''35.66'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -3.00%  '
$ws.Range('D39').Value = '''0.374'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -7.16%  '
$ws.Range('D40').Value = '0.0₃0720'
$ws.Range('E40').Value = '  +8.37%  '
$ws.Range('E41').Value = '  +0.79%  '
$ws.Range('D42').Value = '''2.62'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +6.68%  '
$ws.Range('D43').Value = '2.927.36'
$ws.Range('E43').Value = '  -0.18%  '
$ws.Range('D44').Value = '''0.122'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -3.00%  '
$ws.Range('D45').Value = '''2.98'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +3.82%  '
$ws.Range('E46').Value = '  +0.92%  '
$ws.Range('D47').Value = '''0.0392'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -1.15%  '
$ws.Range('D48').Value = '''2.57'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -4.83%  '
$ws.Range('D49').Value = '''2.99'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +1.73%  '
$ws.Range('D50').Value = '''0.124'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -1.89%  '
$ws.Range('D51').Value = '''135.48'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +2.97%  '
